# Update the "dSF" column (column F) values on Sheet1 to reflect the
# repulled / recalculated data, per the commit:
#   "repull data, push all data, mean calculation"
#
# The workbook is already open; grab the active workbook/worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new dSF (column F) value
$updates = @{
    2  = 6
    3  = -1
    4  = 1
    6  = -2
    7  = 0
    9  = -4
    10 = -1
    11 = 3
    12 = -1
    13 = -5
    14 = -4
    15 = 2
    16 = 6
    17 = -4
    18 = -1
    19 = 9
    20 = 3
    21 = -3
    22 = -2
    23 = -4
    25 = -3
    26 = 1
    27 = -7
    28 = 5
    29 = 1
    30 = 4
    31 = 2
    33 = -1
    35 = 1
    36 = 2
    37 = -3
    38 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
